# Refresh the cryptos price/volume table to the latest scrape values.
# Cells that would otherwise be auto-parsed by Excel as numbers (plain
# "123.45"-style decimals) are explicitly forced to Text format first so
# they keep being stored as strings, matching the sheet's original layout
# (two-part "thousand.decimal" values like "29.951.00" and the
# percentage strings are already text-safe and need no such nudge).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.951.00"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.907.60"
$ws.Range("E3").Value = "  +0.69%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.8027"
$ws.Range("E5").Value = "  +5.94%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "241.31"
$ws.Range("E6").Value = "  +0.58%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.001"
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3150"
$ws.Range("E8").Value = "  +3.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "26.32"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06900"
$ws.Range("E10").Value = "  +1.07%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07992"
$ws.Range("E11").Value = "  +0.18%  "
$ws.Range("D12").Value = "1.915.06"
$ws.Range("E12").Value = "  +1.04%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.7360"
$ws.Range("E13").Value = "  -1.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.186"
$ws.Range("E14").Value = "  -0.44%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.04"
$ws.Range("E15").Value = "  +2.18%  "
$ws.Range("D16").Value = "29.966.67"
$ws.Range("E16").Value = "  +0.26%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "13.96"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "5.867"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.42"
$ws.Range("E19").Value = "  +3.65%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007717"
$ws.Range("E20").Value = "  +0.32%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.002"
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").Value = "2.156.63"
$ws.Range("E22").Value = "  +0.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.002"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.856"
$ws.Range("E24").Value = "  -1.53%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.74"
$ws.Range("E25").Value = "  +1.33%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.206"
$ws.Range("E26").Value = "  -0.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1427"
$ws.Range("E27").Value = "  +10.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.91"
$ws.Range("E28").Value = "  +1.05%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.026"
$ws.Range("E29").Value = "  -0.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").Value = "  +1.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.514"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.304"
$ws.Range("E32").Value = "  +0.56%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.070"
$ws.Range("E33").Value = "  +1.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05505"
$ws.Range("E34").Value = "  +2.67%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.261"
$ws.Range("E35").Value = "  +0.95%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7325"
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.721"
$ws.Range("E37").Value = "  -0.12%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01922"
$ws.Range("E38").Value = "  +0.03%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.792"
$ws.Range("E39").Value = "  +0.90%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.153"
$ws.Range("E40").Value = "  -0.38%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4411"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.20"
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("E43").Value = "  +0.15%  "
$ws.Range("E44").Value = "  +1.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.872"
$ws.Range("E45").Value = "  -2.25%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "100.44"
$ws.Range("E46").Value = "  -0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.552"
$ws.Range("E47").Value = "  -0.24%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "978.36"
$ws.Range("E48").Value = "  +6.41%  "
$ws.Range("B49").Value = "RocketPoolETH"
$ws.Range("C49").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D49").Value = "2.063.95"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.23"
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05965"
$ws.Range("E51").Value = "  +0.04%  "
